$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "amt_104_tk"

# Header row
$ws.Range("A1").Value = 'entryCode'
$ws.Range("B1").Value = 'total'
$ws.Range("C1").Value = 'total.valid'
$ws.Range("D1").Value = 'Worker.ID'
$ws.Range("E1").Value = 'toRate'
$ws.Range("F1").Value = 'Comment'
$ws.Range("G1").Value = 'turker.Index'

# Data rows
# row 2
$ws.Range("A2").Value = '0449e1ec3a904f407b0a7d2f4c9c0c79'
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 'A36D0LIE5AKL2P'
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 184
# row 3
$ws.Range("A3").Value = '0d315e83d7930f9833a05911063d34a9'
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 'A3T111ZNMZ7A5D'
$ws.Range("E3").Value = 0
$ws.Range("E3").Font.Color = 255
$ws.Range("F3").Font.Color = 255
$ws.Range("F3").Value = 'Sound is choppy for video 8'
# row 4
$ws.Range("A4").Value = '0f2c9354f74342a454c4cf9f7e49963f'
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 'A2ZGQUSBB0TMK4'
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = 185
# row 5
$ws.Range("A5").Value = '19dd3ea9be71bfae5cc727811e99d17c'
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = ' AUGR3YWEMVIY7'
$ws.Range("E5").Value = 1
$ws.Range("G5").Value = 186
# row 6
$ws.Range("A6").Value = '1af3743c21839da61d181cbb5803c340'
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 'A1D0JL4MJ7XSGV'
$ws.Range("E6").Value = 0
$ws.Range("E6").Font.Color = 255
$ws.Range("F6").Font.Color = 255
$ws.Range("F6").Value = 'sound is choppy for video 1 to 4'
# row 7
$ws.Range("A7").Value = '1d5ff8280452bb27571b08f5927b56d4'
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 'A330ISTI9O9WTI'
$ws.Range("E7").Value = 0
$ws.Range("E7").Font.Color = 255
$ws.Range("F7").Font.Color = 255
$ws.Range("F7").Value = 'no audio'
# row 8
$ws.Range("A8").Value = '28e1ee5d09ef573b627af22c47ad6f63'
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 'A1ZB2NY0F9QNP0'
$ws.Range("E8").Value = 1
$ws.Range("G8").Value = 187
# row 9
$ws.Range("A9").Value = '2b5326c191e10fcddecf9f7fb478bbc7'
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'A1WR6M74EOTJNY'
$ws.Range("E9").Value = 1
$ws.Range("G9").Value = 188
# row 10
$ws.Range("A10").Value = '2c2da194a6086709c458ebb5097e45e2'
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 'AR72L0JX4D03W'
$ws.Range("E10").Value = 1
$ws.Range("G10").Value = 189
# row 11
$ws.Range("A11").Value = '2e433b6fa7238ec15676dc3269081849'
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 'A2JLSY93R8P8DS'
$ws.Range("E11").Value = 0
$ws.Range("E11").Font.Color = 255
$ws.Range("F11").Font.Color = 255
$ws.Range("F11").Value = 'video quality is poor.'
$ws.Range("A11").NumberFormat = "0.00E+00"
# row 12
$ws.Range("A12").Value = '339bffeb4ebd89494888374822683baa'
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 'A325BIJIG3AK1T'
$ws.Range("E12").Value = 1
$ws.Range("G12").Value = 190
# row 13
$ws.Range("A13").Value = '3970d2c3503b521357dc741561c86220'
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 'A121W3A5FW3MD4'
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 191
# row 14
$ws.Range("A14").Value = '39def480018ef9bb07a1d8664f7e1257'
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 8
$ws.Range("D14").Value = 'A4ZW4GNQ98HV6'
$ws.Range("E14").Value = 1
$ws.Range("G14").Value = 192
# row 15
$ws.Range("A15").Value = '3e084d125ad51c6de592f93d7540c59b'
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 'A3PGUPNMOU5BPW'
$ws.Range("E15").Value = 0
$ws.Range("E15").Font.Color = 255
$ws.Range("F15").Font.Color = 255
$ws.Range("F15").Value = 'no audio'
# row 16
$ws.Range("A16").Value = '41ae951683f8a81bd4c3d7d1e1e57547'
$ws.Range("B16").Value = 8
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 'A3RIFUEQ95MR16'
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = 193
# row 17
$ws.Range("A17").Value = '5522cc12c09f2469f7cc856bc0159aaf'
$ws.Range("B17").Value = 8
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 'A1GOJEDZM2CQTN'
$ws.Range("E17").Value = 1
$ws.Range("G17").Value = 194
# row 18
$ws.Range("A18").Value = '5dd3e271ccd3b8d7f3e3f47985bfada9'
$ws.Range("B18").Value = 8
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 'A3FC1KVPO4RGIN'
$ws.Range("E18").Value = 1
$ws.Range("G18").Value = 195
# row 19
$ws.Range("A19").Value = '5e353316e8373543cd2e041c303bda16'
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 'A9UCDP0W2FVAV'
$ws.Range("E19").Value = 1
$ws.Range("G19").Value = 196
$ws.Range("A19").NumberFormat = "0.00E+00"
# row 20
$ws.Range("A20").Value = '6c516879d637b6746ce52c55af5cf5cc'
$ws.Range("B20").Value = 8
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 'A18C452IX8HCB8'
$ws.Range("E20").Value = 1
$ws.Range("G20").Value = 197
# row 21
$ws.Range("A21").Value = '716367eb320d1ee182e190690a93a8e0'
$ws.Range("B21").Value = 8
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 'A2ISX4NLTUMNPD'
$ws.Range("E21").Value = 1
$ws.Range("G21").Value = 198
# row 22
$ws.Range("A22").Value = '7fa1bd5284f4c08135070397228f9bb9'
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 'A2YV50DPZCX2JN'
$ws.Range("E22").Value = 1
$ws.Range("G22").Value = 199
# row 23
$ws.Range("A23").Value = '8515d632856d816e8e5d3ec96f929877'
$ws.Range("B23").Value = 8
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 'A1YGENVE7OI5JL'
$ws.Range("E23").Value = 0
$ws.Range("E23").Font.Color = 255
$ws.Range("F23").Font.Color = 255
$ws.Range("F23").Value = 'web cam is too old. Just got very small videos'
# row 24
$ws.Range("A24").Value = '85d1efe45a2b7ab2c0f9a592e9673dfa'
$ws.Range("B24").Value = 8
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 'ASX5NVC2MTJ3B'
$ws.Range("E24").Value = 1
$ws.Range("G24").Value = 200
# row 25
$ws.Range("A25").Value = '884cd6903eb684d67ebb4ef3f5e54e22'
$ws.Range("B25").Value = 8
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 'A2UCSSENYIMEIW'
$ws.Range("E25").Value = 1
$ws.Range("G25").Value = 201
# row 26
$ws.Range("A26").Value = '89e9f1b038cf56a572049c20e6fe598a'
$ws.Range("B26").Value = 8
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 'A36A8J9M8GSIYF'
$ws.Range("E26").Value = 1
$ws.Range("G26").Value = 202
# row 27
$ws.Range("A27").Value = '8d411380268b2935eeec299b479df8d5'
$ws.Range("B27").Value = 8
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 'A3TZTSRQJQ97L8'
$ws.Range("E27").Value = 1
$ws.Range("G27").Value = 203
# row 28
$ws.Range("A28").Value = '8e8d7bd42c782f1a84b97b80bec7433f'
$ws.Range("B28").Value = 8
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 'A37EV8RZ82WT8E'
$ws.Range("E28").Value = 0
$ws.Range("E28").Font.Color = 255
$ws.Range("F28").Font.Color = 255
$ws.Range("F28").Value = 'no audio'
# row 29
$ws.Range("A29").Value = '9232a51f2edc9a32d8eeb8f93e6e8a63'
$ws.Range("B29").Value = 8
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 'AAYRVW54BRE8R'
$ws.Range("E29").Value = 1
$ws.Range("G29").Value = 204
# row 30
$ws.Range("A30").Value = '9469fada413f26e0131e366202714c94'
$ws.Range("B30").Value = 8
$ws.Range("C30").Value = 8
$ws.Range("D30").Value = 'A2JQ7V9EWJ51T2'
$ws.Range("E30").Value = 1
$ws.Range("G30").Value = 205
# row 31
$ws.Range("A31").Value = '9f8a8ef97677f66537e0d99bb2478f48'
$ws.Range("B31").Value = 8
$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 'AQCWC087K8YCR'
$ws.Range("E31").Value = 1
$ws.Range("G31").Value = 206
# row 32
$ws.Range("A32").Value = 'a14058177bdb2089e226ef337b8bda36'
$ws.Range("B32").Value = 8
$ws.Range("C32").Value = 8
$ws.Range("D32").Value = 'A29X709ZWO05H2'
$ws.Range("E32").Value = 0
$ws.Range("E32").Font.Color = 255
$ws.Range("F32").Font.Color = 255
$ws.Range("F32").Value = 'partial face'
# row 33
$ws.Range("A33").Value = 'a3f5c6f6bb44cbbb1e5fd97413ca2f50'
$ws.Range("B33").Value = 8
$ws.Range("C33").Value = 8
$ws.Range("D33").Value = 'A3DH2RU1CFTXAP'
$ws.Range("E33").Value = 1
$ws.Range("G33").Value = 207
# row 34
$ws.Range("A34").Value = 'adc04da710e33abdd62c8ec395a15e34'
$ws.Range("B34").Value = 8
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = 'AR5E0ZVWKJA95'
$ws.Range("E34").Value = 1
$ws.Range("G34").Value = 208
# row 35
$ws.Range("A35").Value = 'b1764000277fa055127df47c84b72dab'
$ws.Range("B35").Value = 8
$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 'A3RS7UCO7CQ74R'
$ws.Range("E35").Value = 1
$ws.Range("G35").Value = 209
# row 36
$ws.Range("A36").Value = 'b682b80e30156e2304d4bd92ce8c563e'
$ws.Range("B36").Value = 8
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 'A2EBQ6NGXSXRW0'
$ws.Range("E36").Value = 0
$ws.Range("E36").Font.Color = 255
$ws.Range("F36").Font.Color = 255
$ws.Range("F36").Value = 'no audio'
# row 37
$ws.Range("A37").Value = 'b8b503c7a8bbbecc84c8e448a05cefcb'
$ws.Range("B37").Value = 8
$ws.Range("C37").Value = 8
$ws.Range("D37").Value = 'A1CTOT46Y4W11J'
$ws.Range("E37").Value = 1
$ws.Range("G37").Value = 210
# row 38
$ws.Range("A38").Value = 'c1500ef2b6fe0b3a83a5a598d8c2310f'
$ws.Range("B38").Value = 8
$ws.Range("C38").Value = 8
$ws.Range("D38").Value = 'ARD3X7QES1YF3'
$ws.Range("E38").Value = 1
$ws.Range("G38").Value = 211
# row 39
$ws.Range("A39").Value = 'd4e7a8b9d1e10a53d6cab0d4c50850e9'
$ws.Range("B39").Value = 8
$ws.Range("C39").Value = 8
$ws.Range("D39").Value = 'A1TISWAW29WUGA'
$ws.Range("E39").Value = 1
$ws.Range("G39").Value = 212
# row 40
$ws.Range("A40").Value = 'd93e74eae15c4983f5ffb43b45febf03'
$ws.Range("B40").Value = 8
$ws.Range("C40").Value = 8
$ws.Range("D40").Value = 'A2OR95QZT4H80T'
$ws.Range("E40").Value = 1
$ws.Range("G40").Value = 213
# row 41
$ws.Range("A41").Value = 'd9579a69358f7e16c65e60ba08e14201'
$ws.Range("B41").Value = 8
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = 'A34W8AC87LAAS'
$ws.Range("E41").Value = 1
$ws.Range("G41").Value = 214
# row 42
$ws.Range("A42").Value = 'de56fac608c314d8c1c5370d63df0a7b'
$ws.Range("B42").Value = 8
$ws.Range("C42").Value = 8
$ws.Range("D42").Value = 'A1QG4N21BF61PC'
$ws.Range("E42").Value = 1
$ws.Range("G42").Value = 215
# row 43
$ws.Range("A43").Value = 'de5d1b5a878c07c5e8d3f43cbedcf2c8'
$ws.Range("B43").Value = 8
$ws.Range("C43").Value = 8
$ws.Range("D43").Value = 'A1CJM3ULFBWN1E'
$ws.Range("E43").Value = 1
$ws.Range("G43").Value = 216
# row 44
$ws.Range("A44").Value = 'e256153fbe35829495511810e3dc96ff'
$ws.Range("B44").Value = 8
$ws.Range("C44").Value = 8
$ws.Range("D44").Value = 'A5NKVJX6QEXEY'
$ws.Range("E44").Value = 1
$ws.Range("G44").Value = 217
# row 45
$ws.Range("A45").Value = 'effc4b0420cb5df36b01e4954372d2cd'
$ws.Range("B45").Value = 8
$ws.Range("C45").Value = 8
$ws.Range("D45").Value = 'A2VBSFSJXLZZ7A'
$ws.Range("E45").Value = 1
$ws.Range("G45").Value = 218
# row 46
$ws.Range("A46").Value = 'f5148d89f7857408575204c24bb12714'
$ws.Range("B46").Value = 8
$ws.Range("C46").Value = 8
$ws.Range("D46").Value = 'ALF9AAZGQP4K5'
$ws.Range("E46").Value = 1
$ws.Range("G46").Value = 219
# row 47
$ws.Range("A47").Value = 'ffb621477e000af8f66717bcd6278482'
$ws.Range("B47").Value = 8
$ws.Range("C47").Value = 8
$ws.Range("D47").Value = 'A60BX1JSJRYAX'
$ws.Range("E47").Value = 1
$ws.Range("G47").Value = 220
# row 48
$ws.Range("A48").Value = 'e3809ec34cf9f345ed4e59ab033f9db3'
$ws.Range("B48").Value = 8
$ws.Range("C48").Value = 7
$ws.Range("D48").Value = 'A260QC80Q9VX7Z'
$ws.Range("E48").Value = 0
$ws.Range("E48").Font.Color = 255
$ws.Range("F48").Font.Color = 255
# row 49
$ws.Range("A49").Value = '7d5799f278b55122c8de5542180d43ab'
$ws.Range("B49").Value = 8
$ws.Range("C49").Value = 6
$ws.Range("D49").Value = 'A194E79BDFEZIU'
$ws.Range("E49").Value = 0
$ws.Range("E49").Font.Color = 255
$ws.Range("F49").Font.Color = 255
# row 50
$ws.Range("A50").Value = '8a530de58927745c927fc8ff752d0ee6'
$ws.Range("B50").Value = 8
$ws.Range("C50").Value = 6
$ws.Range("D50").Value = 'A37H5799TLDXUT'
$ws.Range("E50").Value = 0
$ws.Range("E50").Font.Color = 255
$ws.Range("F50").Font.Color = 255
# row 51
$ws.Range("A51").Value = 'de82a643ed658b40447cfcb85249956e'
$ws.Range("B51").Value = 8
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 'A2FV527DH1S2Y7'
$ws.Range("E51").Value = 0
$ws.Range("E51").Font.Color = 255
$ws.Range("F51").Font.Color = 255
# row 52
$ws.Range("A52").Value = '85da6399de4f8791b7e5f4a36d81fbc7'
$ws.Range("B52").Value = 8
$ws.Range("C52").Value = 3
$ws.Range("D52").Value = 'A2WCV7ULUPKVNQ'
$ws.Range("E52").Value = 0
$ws.Range("E52").Font.Color = 255
$ws.Range("F52").Font.Color = 255
# row 53
$ws.Range("A53").Value = '3f10bb3b932ba8b6f95fffd18b69dbc1'
$ws.Range("B53").Value = 8
$ws.Range("C53").Value = 2
$ws.Range("D53").Value = 'A1S8DYWNS59XWB'
$ws.Range("E53").Value = 0
$ws.Range("E53").Font.Color = 255
$ws.Range("F53").Font.Color = 255

# Column widths
$ws.Columns.Item(1).ColumnWidth = 33.666666666666664
$ws.Columns.Item(4).ColumnWidth = 22.833333333333336
$ws.Columns.Item(6).ColumnWidth = 39.166666666666664
$ws.Columns.Item(7).ColumnWidth = 9.5

# Selection
$null = $ws.Range("E48:F53").Select()
